$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 1: "Method recover() implements error recovery by skipping tokens
# until it finds one in the follow set of the nonterminal defined by the rule."
# Split run 3 so the sentence reads "...until it finds one whose symbol is in
# the follow set of the nonterminal defined by the rule."
$para1 = $tr.Paragraphs(1, 1)
$run3 = $para1.Runs(3, 1)
$run3.Text = " implements error recovery by skipping tokens until it finds one "
$run3.InsertAfter("whose symbol is in the follow set of the nonterminal defined by the rule.") | Out-Null
